## Apply the "Adding small models to test" edit:
##  - sheet "constraints": shift the old running-total column H into a new
##    column G (with its original formula), rewrite column H with a new set
##    of per-row running totals, and rewrite column I with a new set of
##    running totals that occasionally branches off column H (submodel
##    linking). Highlight a few cells (I4, I11, I16, I18) with the green
##    fill already used elsewhere in the sheet, and mark rows 15-18's
##    label cells (A15:A18) with the yellow fill used by the other rows.
##  - sheet "variables": add a new column F of running totals that chains
##    off the "constraints" sheet's totals, and recolor A2:A9 from yellow
##    to the green fill already used by A10:A11.

$wb = $excel.ActiveWorkbook

$wsC = $wb.Worksheets.Item("constraints")
$wsV = $wb.Worksheets.Item("variables")

# ---------------------------------------------------------------------
# Sheet "constraints": columns G, H, I for rows 2-18
# ---------------------------------------------------------------------

# Column G: running total that continues the sheet's original column-H
# pattern (base case 48+25+C2, then C{r}+G{r-1}).
$wsC.Range("G2").Formula = "=48+25+C2"
for ($r = 3; $r -le 18; $r++) {
    $wsC.Range("G$r").Formula = "=C$r+G" + ($r - 1)
}

# Column H: new set of per-row formulas (no longer a single shared
# formula - each row has its own explicit relationship).
$wsC.Range("H2").Formula  = "=C2"
$wsC.Range("H3").Formula  = "=H2+C3"
$wsC.Range("H4").Formula  = "=C4+H3"
$wsC.Range("H5").Formula  = "=48+C5"
$wsC.Range("H6").Formula  = "=C6+H5"
$wsC.Range("H7").Formula  = "=C7+H6"
$wsC.Range("H8").Formula  = "=C8+H7"
$wsC.Range("H9").Formula  = "=C9+H8"
$wsC.Range("H10").Formula = "=C10+H9"
$wsC.Range("H11").Formula = "=C11+H10"
$wsC.Range("H12").Formula = "=25+C12"
$wsC.Range("H13").Formula = "=H12+C13"
$wsC.Range("H14").Formula = "=H13+C14"
$wsC.Range("H15").Formula = "=H14+C15"
$wsC.Range("H16").Formula = "=H15+C16"
$wsC.Range("H17").Formula = "=C17"
$wsC.Range("H18").Formula = "=H17+C18"

# Column I: new set of per-row formulas; some rows link back into column
# H of a different row (the "submodel" cross reference called out in the
# commit message).
$wsC.Range("I2").Formula  = "=H4+C2"
$wsC.Range("I3").Formula  = "=I2+C3"
$wsC.Range("I4").Formula  = "=I3+C4"
$wsC.Range("I5").Formula  = "=H11+C5"
$wsC.Range("I6").Formula  = "=I5+C6"
$wsC.Range("I7").Formula  = "=I6+C7"
$wsC.Range("I8").Formula  = "=I7+C8"
$wsC.Range("I9").Formula  = "=I8+C9"
$wsC.Range("I10").Formula = "=I9+C10"
$wsC.Range("I11").Formula = "=I10+C11"
$wsC.Range("I12").Formula = "=H16+C12"
$wsC.Range("I13").Formula = "=I12+C13"
$wsC.Range("I14").Formula = "=I13+C14"
$wsC.Range("I15").Formula = "=I14+C15"
$wsC.Range("I16").Formula = "=I15+C16"
$wsC.Range("I17").Formula = "=H18+C17"
$wsC.Range("I18").Formula = "=I17+C18"

# Highlight the "submodel handoff" cells with the green fill already used
# by the other shaded cells on the sheet (same fill as row 19's label).
$greenColor = $wsC.Range("A19").Interior.Color
foreach ($addr in @("I4", "I11", "I16", "I18")) {
    $wsC.Range($addr).Interior.Color = $greenColor
}

# Mark the labels for rows 15-18 with the yellow fill used by rows 2-14's
# labels.
$yellowColor = $wsC.Range("A14").Interior.Color
foreach ($addr in @("A15", "A16", "A17", "A18")) {
    $wsC.Range($addr).Interior.Color = $yellowColor
}

# Move the sheet's selection cursor.
[void]$wsC.Range("L14").Select()

# ---------------------------------------------------------------------
# Sheet "variables": column F for rows 2-9, chaining off constraints!I
# ---------------------------------------------------------------------

$wsV.Range("F2").Formula = "=4440+B2"
$wsV.Range("F3").Formula = "=F2+B3"
for ($r = 4; $r -le 9; $r++) {
    $wsV.Range("F$r").Formula = "=F" + ($r - 1) + "+B$r"
}

# Recolor A2:A9 from yellow to the green fill already used by A10:A11.
$greenColorV = $wsV.Range("A10").Interior.Color
$wsV.Range("A2:A9").Interior.Color = $greenColorV
